# issue #5: add legislator_id, name, date into dataframe
#
# The 股票 (stocks) sheet gains three new trailing columns -- date,
# legislator_name, legislator_id -- populated with the same value on every
# data row (this property-disclosure row came from one filing, so every
# row shares the filer's name/id and the filing date).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$filingDate   = "2012-04-24"
$legislator   = "楊玉欣"
$legislatorId = 1757

# data rows are 2..9 (row 1 is the header)
$lastRow = $ws.UsedRange.Rows.Count

# --- headers (H1:J1), matching the look of the existing B1:G1 headers ---
$ws.Range("H1:J1").Font.Bold = $true
$ws.Range("H1:J1").Borders.LineStyle = 1
$ws.Range("H1:J1").HorizontalAlignment = -4108
$ws.Range("H1:J1").VerticalAlignment = -4160

$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# Format the date column as text first so the "YYYY-MM-DD" string is kept
# literally instead of being auto-converted into a date serial number.
$ws.Range("H2:H" + $lastRow).NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value  = $filingDate
    $ws.Cells.Item($r, 9).Value  = $legislator
    $ws.Cells.Item($r, 10).Value = $legislatorId
}
